$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Gender" column right after "Subject No." (col A) ---
# This shifts the old B..M (height .. distance/right ankle) over to C..N,
# matching the sharedStrings / cell-ref reshuffle in the target file.
$ws.Columns.Item(2).Insert()

# Approximate the inserted column's width to match column A's custom width
# (Excel's native "Insert Column" carries the left neighbour's width along).
$ws.Columns.Item(2).ColumnWidth = 12.83

# Header for the new column
$ws.Range("B1").Value = "Gender"

# Row 3 ("subject 2") carries a row-level custom format (s=4 on the row,
# s=3 on its cells); the target file still wants the new Gender cell in
# that row to use the plain/default cell style like every other Gender
# cell (s=1), so copy that formatting over from B2 before filling values.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

# Gender values, one per subject row (2-15)
$genders = @("Male","Male","Male","Male","Female","Male","Female","Female","Male","Male","Female","Female","Male","Male")
for ($i = 0; $i -lt $genders.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $genders[$i]
}

# Match the saved selection/active cell from the target workbook
$ws.Range("D18").Select()
